# v1.0.19 Updating for different macOS devices. Addressing reviewer comments.
#
# The underlying Cox-regression figures (Beta, s.e.m., HR, CIs, Z, P) were
# recomputed on a different machine/OS; only the last one or two
# significant digits of a subset of the numeric results cells moved
# (classic cross-platform libm last-bit rounding). Re-apply the refreshed
# values onto the "Results" sheet, cell by cell, and nudge the sheet's
# default column width metadata to match the value macOS Excel writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Sheet-level formatting nudge that shows up when the file is re-saved on
# a macOS Excel install (sheetFormatPr baseColWidth="10").
$ws.StandardWidth = 10

$updates = @{
    "E2" = 0.405328598683507
    "F2" = 2.94017599429189
    "G2" = 1.32846665762654
    "I2" = 2.66072871508649
    "J2" = 0.00779717513921154

    "D4" = 0.840815927943706
    "G4" = 1.00243924875164
    "H4" = 5.36124153151077
    "J4" = 0.0493337791244438

    "I5" = 0.335138719090466

    "F6" = 7.83908191961647
    "G6" = 0.910282757444315
    "H6" = 67.507820883026
    "J6" = 0.0608710144622736

    "D7" = -0.777360451280068
    "H7" = 5.08135055567079
    "I7" = -0.634068286148468
    "J7" = 0.526036258577317

    "D8" = 65.4956835568229
    "E8" = 3024.67962969327
    "F8" = [double]"2.78236407491737E+28"
    "I8" = 0.021653758934947
    "J8" = 0.98272415013828

    "D9" = 22.7617208196112
    "E9" = 51224.636929086
    "F9" = 7678736248.66679
    "I9" = 0.000444351042470479
    "J9" = 0.999645459175303
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value2 = $updates[$addr]
}
